$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "(Datum)" -> "Datum" (3 occurrences, each inside a highlighted run)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("(Datum)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Datum", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Frame1 drawing: grow the anchor extent a touch (198.15pt x 141.45pt ->
#    198.2pt x 141.5pt), matching the new <wp:extent>.
# ---------------------------------------------------------------------------
$frame = $d.Shapes.Item("Frame1")
$frame.Width = 198.2
$frame.Height = 141.5

# ---------------------------------------------------------------------------
# 3) overflowPunct: false -> true, i.e. HangingPunctuation off -> on, on the
#    Normal style and a handful of other paragraph styles.
# ---------------------------------------------------------------------------
$overflowPunctStyles = @("Normal", "TOC 6", "TOC 7", "TOC 8", "TOC 9", "Wappen")
foreach ($styleName in $overflowPunctStyles) {
    $style = $d.Styles.Item($styleName)
    $style.ParagraphFormat.HangingPunctuation = $true
}

# ---------------------------------------------------------------------------
# 4) New character styles ListLabel63 .. ListLabel71 (mirrors the existing
#    ListLabelNN character styles already in the template).
# ---------------------------------------------------------------------------
$newListLabels = @(
    @{Id = "ListLabel63"; Font = "Arial"},
    @{Id = "ListLabel64"; Font = "Courier New"},
    @{Id = "ListLabel65"; Font = "Wingdings"},
    @{Id = "ListLabel66"; Font = "Symbol"},
    @{Id = "ListLabel67"; Font = "Courier New"},
    @{Id = "ListLabel68"; Font = "Wingdings"},
    @{Id = "ListLabel69"; Font = "Symbol"},
    @{Id = "ListLabel70"; Font = "Courier New"},
    @{Id = "ListLabel71"; Font = "Wingdings"}
)

foreach ($def in $newListLabels) {
    $style = $d.Styles.Add($def.Id, 2)
    $style.NameLocal = ($def.Id -replace "ListLabel", "ListLabel ")
    $style.QuickStyle = $true
    $style.Font.NameBi = $def.Font
}
